$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.924.97'
$ws.Range('E2').Value = '  -2.03%  '
$ws.Range('D3').Value = '1.901.05'
$ws.Range('E3').Value = '  -4.02%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.00'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4580'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3809'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07708'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9736'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.96'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.14%  '
$ws.Range('D12').Value = '1.896.38'
$ws.Range('E12').Value = '  -5.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.912'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.643'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07054'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '83.63'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009463'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -5.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.60'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.87%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '28.898.56'
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.272'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.13%  '
$ws.Range('E23').Value = '  -3.06%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.095'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.08'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.01'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.17%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.602'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.88%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '117.35'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.828'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.15%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09239'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.8560'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.01%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.072'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.235'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -6.79%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.999'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.28%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.05653'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.140'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.27%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02028'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.41%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5472'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.65%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.377'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.40%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1750'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.06%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.248'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.48%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.753'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.26%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5139'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.32%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.11'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.63%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.06812'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.063'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.78%  '
$ws.Range('B48').Value = 'PEPE'
$ws.Range('C48').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000002583'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -15.04%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '109.98'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.77%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.765'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.55%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.15%  '
